{"js": "// Update the date paragraph (first paragraph in the body) and the\n// multiplication-table answers to match the edited document.\n//\n// Expected (\"before\") values are embedded and checked defensively: a\n// cell is only updated when its current text still matches the known\n// \"before\" value, so the script is a no-op (per cell) if the document\n// has already been edited or differs from what we expect.\n\nconst GRID_OLD = [[\"10\u00d797=970\", \"78\u00d720=1560\", \"64\u00d730=1920\", \"60\u00d783=4980\", \"84\u00d764=5376\"], [\"90\u00d799=8910\", \"87\u00d791=7917\", \"40\u00d743=1720\", \"35\u00d717=595\", \"92\u00d776=6992\"], [\"91\u00d784=7644\", \"28\u00d718=504\", \"29\u00d719=551\", \"15\u00d798=1470\", \"84\u00d762=5208\"], [\"69\u00d717=1173\", \"23\u00d793=2139\", \"99\u00d757=5643\", \"14\u00d781=1134\", \"51\u00d757=2907\"], [\"90\u00d774=6660\", \"58\u00d797=5626\", \"15\u00d751=765\", \"69\u00d781=5589\", \"57\u00d767=3819\"], [\"95\u00d752=4940\", \"24\u00d751=1224\", \"61\u00d750=3050\", \"80\u00d763=5040\", \"53\u00d734=1802\"], [\"47\u00d722=1034\", \"30\u00d736=1080\", \"53\u00d797=5141\", \"10\u00d716=160\", \"86\u00d710=860\"], [\"99\u00d759=5841\", \"65\u00d770=4550\", \"25\u00d796=2400\", \"16\u00d767=1072\", \"14\u00d779=1106\"], [\"32\u00d787=2784\", \"14\u00d751=714\", \"92\u00d790=8280\", \"10\u00d722=220\", \"14\u00d735=490\"], [\"65\u00d756=3640\", \"38\u00d765=2470\", \"75\u00d781=6075\", \"29\u00d757=1653\", \"99\u00d797=9603\"], [\"53\u00d722=1166\", \"41\u00d717=697\", \"95\u00d788=8360\", \"10\u00d720=200\", \"39\u00d754=2106\"], [\"77\u00d715=1155\", \"95\u00d763=5985\", \"61\u00d729=1769\", \"14\u00d734=476\", \"60\u00d796=5760\"], [\"55\u00d758=3190\", \"39\u00d732=1248\", \"35\u00d731=1085\", \"90\u00d766=5940\", \"31\u00d720=620\"], [\"13\u00d725=325\", \"23\u00d744=1012\", \"23\u00d769=1587\", \"43\u00d785=3655\", \"22\u00d788=1936\"], [\"67\u00d739=2613\", \"74\u00d784=6216\", \"86\u00d735=3010\", \"99\u00d793=9207\", \"52\u00d761=3172\"], [\"27\u00d718=486\", \"65\u00d798=6370\", \"85\u00d758=4930\", \"92\u00d754=4968\", \"99\u00d719=1881\"], [\"65\u00d783=5395\", \"88\u00d716=1408\", \"49\u00d765=3185\", \"93\u00d754=5022\", \"23\u00d756=1288\"], [\"35\u00d783=2905\", \"47\u00d788=4136\", \"57\u00d757=3249\", \"58\u00d740=2320\", \"87\u00d777=6699\"], [\"21\u00d758=1218\", \"42\u00d766=2772\", \"39\u00d718=702\", \"73\u00d793=6789\", \"98\u00d763=6174\"], [\"28\u00d741=1148\", \"22\u00d723=506\", \"58\u00d753=3074\", \"28\u00d765=1820\", \"64\u00d728=1792\"]];\nconst GRID_NEW = [[\"58\u00d777=4466\", \"18\u00d770=1260\", \"50\u00d798=4900\", \"40\u00d799=3960\", \"92\u00d784=7728\"], [\"94\u00d741=3854\", \"63\u00d772=4536\", \"85\u00d788=7480\", \"39\u00d720=780\", \"59\u00d731=1829\"], [\"52\u00d791=4732\", \"53\u00d795=5035\", \"34\u00d777=2618\", \"38\u00d794=3572\", \"44\u00d791=4004\"], [\"87\u00d738=3306\", \"37\u00d789=3293\", \"57\u00d746=2622\", \"84\u00d711=924\", \"37\u00d767=2479\"], [\"43\u00d723=989\", \"86\u00d734=2924\", \"94\u00d720=1880\", \"34\u00d754=1836\", \"62\u00d730=1860\"], [\"69\u00d724=1656\", \"100\u00d710=1000\", \"91\u00d787=7917\", \"42\u00d742=1764\", \"46\u00d773=3358\"], [\"48\u00d711=528\", \"17\u00d7100=1700\", \"70\u00d733=2310\", \"74\u00d762=4588\", \"23\u00d785=1955\"], [\"93\u00d773=6789\", \"79\u00d747=3713\", \"33\u00d728=924\", \"41\u00d764=2624\", \"80\u00d738=3040\"], [\"58\u00d742=2436\", \"71\u00d737=2627\", \"89\u00d720=1780\", \"15\u00d753=795\", \"31\u00d714=434\"], [\"99\u00d745=4455\", \"59\u00d752=3068\", \"15\u00d717=255\", \"43\u00d783=3569\", \"41\u00d718=738\"], [\"85\u00d735=2975\", \"77\u00d729=2233\", \"54\u00d723=1242\", \"45\u00d734=1530\", \"47\u00d767=3149\"], [\"39\u00d766=2574\", \"55\u00d764=3520\", \"95\u00d774=7030\", \"90\u00d720=1800\", \"83\u00d731=2573\"], [\"72\u00d729=2088\", \"85\u00d794=7990\", \"86\u00d732=2752\", \"36\u00d757=2052\", \"65\u00d747=3055\"], [\"13\u00d725=325\", \"23\u00d758=1334\", \"17\u00d771=1207\", \"29\u00d712=348\", \"30\u00d772=2160\"], [\"59\u00d796=5664\", \"88\u00d769=6072\", \"27\u00d780=2160\", \"49\u00d768=3332\", \"25\u00d722=550\"], [\"14\u00d751=714\", \"68\u00d788=5984\", \"27\u00d776=2052\", \"22\u00d790=1980\", \"61\u00d755=3355\"], [\"59\u00d741=2419\", \"75\u00d722=1650\", \"59\u00d714=826\", \"63\u00d716=1008\", \"11\u00d786=946\"], [\"56\u00d765=3640\", \"54\u00d721=1134\", \"48\u00d719=912\", \"70\u00d779=5530\", \"34\u00d772=2448\"], [\"47\u00d725=1175\", \"32\u00d721=672\", \"54\u00d724=1296\", \"56\u00d724=1344\", \"95\u00d771=6745\"], [\"12\u00d711=132\", \"16\u00d763=1008\", \"71\u00d759=4189\", \"65\u00d771=4615\", \"81\u00d736=2916\"]];\nconst DATE_OLD = \"2023-04-13 Thursday\";\nconst DATE_NEW = \"2023-04-14 Friday\";\n\nconst body = context.document.body;\n\n// --- 1) Update the date heading paragraph -------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\n\nif (dateParagraph.text.trim() === DATE_OLD) {\n  dateParagraph.insertText(DATE_NEW, \"Replace\");\n}\n\n// --- 2) Update the multiplication-table cells ----------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst current = table.values;\nconst rows = GRID_OLD.length;\nconst cols = GRID_OLD[0].length;\n\nfor (let r = 0; r < rows; r++) {\n  for (let c = 0; c < cols; c++) {\n    const expectedOld = GRID_OLD[r][c];\n    const newVal = GRID_NEW[r][c];\n    if (newVal === expectedOld) continue; // unchanged cell, skip\n    if (current[r][c] === expectedOld) {\n      current[r][c] = newVal;\n    }\n  }\n}\n\ntable.values = current;\nawait context.sync();\n", "ps1": "# Update the date paragraph (first paragraph in the body) and the\n# multiplication-table answers to match the edited document.\n#\n# Expected (\"before\") values are checked defensively: a cell/paragraph\n# is only updated when its current text still matches the known \"before\"\n# value, so the script is a no-op where the document differs from what\n# is expected.\n\n$d = $word.ActiveDocument\n\n# --- 1) Update the date heading paragraph --------------------------------\n$dateOld = '2023-04-13 Thursday'\n$dateNew = '2023-04-14 Friday'\n$p1 = $d.Paragraphs.Item(1)\n$p1Text = $p1.Range.Text.TrimEnd([char]13, [char]7)\nif ($p1Text -eq $dateOld) {\n    $p1.Range.Text = $dateNew\n}\n\n# --- 2) Update the multiplication-table cells -----------------------------\n$oldValues = @(\n    @('10\u00d797=970', '78\u00d720=1560', '64\u00d730=1920', '60\u00d783=4980', '84\u00d764=5376'),\n    @('90\u00d799=8910', '87\u00d791=7917', '40\u00d743=1720', '35\u00d717=595', '92\u00d776=6992'),\n    @('91\u00d784=7644', '28\u00d718=504', '29\u00d719=551', '15\u00d798=1470', '84\u00d762=5208'),\n    @('69\u00d717=1173', '23\u00d793=2139', '99\u00d757=5643', '14\u00d781=1134', '51\u00d757=2907'),\n    @('90\u00d774=6660', '58\u00d797=5626', '15\u00d751=765', '69\u00d781=5589', '57\u00d767=3819'),\n    @('95\u00d752=4940', '24\u00d751=1224', '61\u00d750=3050', '80\u00d763=5040', '53\u00d734=1802'),\n    @('47\u00d722=1034', '30\u00d736=1080', '53\u00d797=5141', '10\u00d716=160', '86\u00d710=860'),\n    @('99\u00d759=5841', '65\u00d770=4550', '25\u00d796=2400', '16\u00d767=1072', '14\u00d779=1106'),\n    @('32\u00d787=2784', '14\u00d751=714', '92\u00d790=8280', '10\u00d722=220', '14\u00d735=490'),\n    @('65\u00d756=3640', '38\u00d765=2470', '75\u00d781=6075', '29\u00d757=1653', '99\u00d797=9603'),\n    @('53\u00d722=1166', '41\u00d717=697', '95\u00d788=8360', '10\u00d720=200', '39\u00d754=2106'),\n    @('77\u00d715=1155', '95\u00d763=5985', '61\u00d729=1769', '14\u00d734=476', '60\u00d796=5760'),\n    @('55\u00d758=3190', '39\u00d732=1248', '35\u00d731=1085', '90\u00d766=5940', '31\u00d720=620'),\n    @('13\u00d725=325', '23\u00d744=1012', '23\u00d769=1587', '43\u00d785=3655', '22\u00d788=1936'),\n    @('67\u00d739=2613', '74\u00d784=6216', '86\u00d735=3010', '99\u00d793=9207', '52\u00d761=3172'),\n    @('27\u00d718=486', '65\u00d798=6370', '85\u00d758=4930', '92\u00d754=4968', '99\u00d719=1881'),\n    @('65\u00d783=5395', '88\u00d716=1408', '49\u00d765=3185', '93\u00d754=5022', '23\u00d756=1288'),\n    @('35\u00d783=2905', '47\u00d788=4136', '57\u00d757=3249', '58\u00d740=2320', '87\u00d777=6699'),\n    @('21\u00d758=1218', '42\u00d766=2772', '39\u00d718=702', '73\u00d793=6789', '98\u00d763=6174'),\n    @('28\u00d741=1148', '22\u00d723=506', '58\u00d753=3074', '28\u00d765=1820', '64\u00d728=1792')\n)\n$newValues = @(\n    @('58\u00d777=4466', '18\u00d770=1260', '50\u00d798=4900', '40\u00d799=3960', '92\u00d784=7728'),\n    @('94\u00d741=3854', '63\u00d772=4536', '85\u00d788=7480', '39\u00d720=780', '59\u00d731=1829'),\n    @('52\u00d791=4732', '53\u00d795=5035', '34\u00d777=2618', '38\u00d794=3572', '44\u00d791=4004'),\n    @('87\u00d738=3306', '37\u00d789=3293', '57\u00d746=2622', '84\u00d711=924', '37\u00d767=2479'),\n    @('43\u00d723=989', '86\u00d734=2924', '94\u00d720=1880', '34\u00d754=1836', '62\u00d730=1860'),\n    @('69\u00d724=1656', '100\u00d710=1000', '91\u00d787=7917', '42\u00d742=1764', '46\u00d773=3358'),\n    @('48\u00d711=528', '17\u00d7100=1700', '70\u00d733=2310', '74\u00d762=4588', '23\u00d785=1955'),\n    @('93\u00d773=6789', '79\u00d747=3713', '33\u00d728=924', '41\u00d764=2624', '80\u00d738=3040'),\n    @('58\u00d742=2436', '71\u00d737=2627', '89\u00d720=1780', '15\u00d753=795', '31\u00d714=434'),\n    @('99\u00d745=4455', '59\u00d752=3068', '15\u00d717=255', '43\u00d783=3569', '41\u00d718=738'),\n    @('85\u00d735=2975', '77\u00d729=2233', '54\u00d723=1242', '45\u00d734=1530', '47\u00d767=3149'),\n    @('39\u00d766=2574', '55\u00d764=3520', '95\u00d774=7030', '90\u00d720=1800', '83\u00d731=2573'),\n    @('72\u00d729=2088', '85\u00d794=7990', '86\u00d732=2752', '36\u00d757=2052', '65\u00d747=3055'),\n    @('13\u00d725=325', '23\u00d758=1334', '17\u00d771=1207', '29\u00d712=348', '30\u00d772=2160'),\n    @('59\u00d796=5664', '88\u00d769=6072', '27\u00d780=2160', '49\u00d768=3332', '25\u00d722=550'),\n    @('14\u00d751=714', '68\u00d788=5984', '27\u00d776=2052', '22\u00d790=1980', '61\u00d755=3355'),\n    @('59\u00d741=2419', '75\u00d722=1650', '59\u00d714=826', '63\u00d716=1008', '11\u00d786=946'),\n    @('56\u00d765=3640', '54\u00d721=1134', '48\u00d719=912', '70\u00d779=5530', '34\u00d772=2448'),\n    @('47\u00d725=1175', '32\u00d721=672', '54\u00d724=1296', '56\u00d724=1344', '95\u00d771=6745'),\n    @('12\u00d711=132', '16\u00d763=1008', '71\u00d759=4189', '65\u00d771=4615', '81\u00d736=2916')\n)\n\n$t = $d.Tables.Item(1)\n$rowCount = 20\n$colCount = 5\n\nfor ($r = 0; $r -lt $rowCount; $r++) {\n    for ($c = 0; $c -lt $colCount; $c++) {\n        $expectedOld = $oldValues[$r][$c]\n        $newVal = $newValues[$r][$c]\n        if ($newVal -eq $expectedOld) { continue }\n        $cell = $t.Cell($r + 1, $c + 1)\n        $cellText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        if ($cellText -eq $expectedOld) {\n            $cell.Range.Text = $newVal\n        }\n    }\n}\n\n"}
